# Update gh-pages to output generated at 456a3b4
# Bumps "想去人数" (F column) counts on the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -----------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 264
$ws1.Range("F3").Value = 70
$ws1.Range("F4").Value = 14
$ws1.Range("F5").Value = 6646
$ws1.Range("F6").Value = 5429
$ws1.Range("F7").Value = 450
$ws1.Range("F11").Value = 235
$ws1.Range("F12").Value = 65

# --- Sheet "全部类型" --------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 264
$ws4.Range("F3").Value = 70
$ws4.Range("F4").Value = 14
$ws4.Range("F5").Value = 6646
$ws4.Range("F6").Value = 5429
$ws4.Range("F7").Value = 450
$ws4.Range("F11").Value = 235
$ws4.Range("F14").Value = 65
